$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Fill new row 2 with the new "Personalnummer" attribute
$ws.Cells.Item(2, 1).Value = "Personalnummer"
$ws.Cells.Item(2, 2).Value = "M100001"

# Apply the label formatting (used by A3, the former row 2) to both new cells
$ws.Range("A3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the selection to match the recorded state
$ws.Range("B7").Select()
